$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "endOfTestData" marker rows get a new yellow highlight fill applied
# (border/other formatting stays the same, only the fill color changes).
$ws.Range("A4").Interior.Color = 65535
$ws.Range("A8").Interior.Color = 65535
$ws.Range("A14").Interior.Color = 65535

# New test-data block for "validateCreateCustomerAPI" is appended after the
# existing "validateCreateUserAPIWithValidData" block (rows 15-18), mirroring
# the other blocks' layout: a header row, a couple of blank rows, then the
# "endOfTestData" marker row.

# Row 15: header row, same look as the other section headers (rows 1/5/9).
$ws.Range("A5:D5").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)
$ws.Range("A15").Value = "validateCreateCustomerAPI"

# Rows 16-17: blank data rows, same plain bordered look as other data rows.
$ws.Range("A6:D7").Copy()
$ws.Range("A16:D17").PasteSpecial(-4122)
$ws.Range("A16:D17").ClearContents()

# Row 18: closing "endOfTestData" marker, highlighted like the other markers.
$ws.Range("A14:D14").Copy()
$ws.Range("A18:D18").PasteSpecial(-4122)
$ws.Range("A18").Value = "endOfTestData"

$excel.CutCopyMode = $false

# Move the active selection to A13.
[void]$ws.Range("A13").Select()
